# [Composer CMS v2.0.beta5 :: Thu, 02 Apr 2015 10:21:53 -0700]
#
# Bumps the version/date line shown on the title page of the generated
# "Composer CMS" user guide from "v2.0.beta4 (2015-04-01)" to
# "v2.0.beta5 (2015-04-02)". Both the version token and the date token
# live in their own runs, so each is matched/replaced independently
# (wdFindContinue keeps searching the same story for the second hit).

$d = $word.ActiveDocument

$wdFindContinue  = 1
$wdReplaceOne    = 1

# "v2.0.beta4" -> "v2.0.beta5"
$d.Content.Find.Execute(
    "v2.0.beta4", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "v2.0.beta5", $wdReplaceOne)

# "(2015-04-01)" -> "(2015-04-02)"
$d.Content.Find.Execute(
    "(2015-04-01)", $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false, "(2015-04-02)", $wdReplaceOne)
